$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list with latest scraped prices / 1h volume changes.
# Note: a leading apostrophe is used for Price values that would otherwise be
# auto-parsed by Excel as numbers, so they stay plain text like the source data
# (e.g. "7.50" keeps its trailing zero instead of becoming the number 7.5).
$ws.Range("D2").Value = '68.348.17'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '2.640.45'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''597.12'
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("D6").Value = '''154.42'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -0.68%  '
$ws.Range("D9").Value = '2.640.44'
$ws.Range("E9").Value = '  +0.28%  '
$ws.Range("E10").Value = '  +7.57%  '
$ws.Range("E11").Value = '  -0.52%  '
$ws.Range("E12").Value = '  +1.03%  '
$ws.Range("E13").Value = '  +1.91%  '
$ws.Range("D14").Value = '''28.15'
$ws.Range("E14").Value = '  +1.63%  '
$ws.Range("E15").Value = '  +2.61%  '
$ws.Range("D16").Value = '3.122.27'
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").Value = '68.327.05'
$ws.Range("E17").Value = '  +0.68%  '
$ws.Range("D18").Value = '2.639.28'
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").Value = '''11.41'
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").Value = '''364.04'
$ws.Range("E20").Value = '  -3.03%  '
$ws.Range("D21").Value = '''7.50'
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("E22").Value = '  +3.23%  '
$ws.Range("E23").Value = '  +1.51%  '
$ws.Range("E24").Value = '  +0.84%  '
$ws.Range("D25").Value = '''74.64'
$ws.Range("E25").Value = '  +3.14%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.11%  '
$ws.Range("E27").Value = '  -1.31%  '
$ws.Range("E28").Value = '  +1.61%  '
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").Value = '''572.48'
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("D32").Value = '''8.16'
$ws.Range("E32").Value = '  +3.84%  '
$ws.Range("E33").Value = '  +1.21%  '
$ws.Range("E34").Value = '  +1.18%  '
$ws.Range("E35").Value = '  +3.26%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").Value = '''1.59'
$ws.Range("E37").Value = '  +4.80%  '
$ws.Range("D38").Value = '''161.04'
$ws.Range("E38").Value = '  +1.70%  '
$ws.Range("D39").Value = '''19.39'
$ws.Range("E39").Value = '  +1.15%  '
$ws.Range("D41").Value = '''1.89'
$ws.Range("E41").Value = '  -0.78%  '
$ws.Range("D42").Value = '''5.41'
$ws.Range("E42").Value = '  +1.22%  '
$ws.Range("D43").Value = '0.0₆0337'
$ws.Range("E43").Value = '  +4.79%  '
$ws.Range("E44").Value = '  +1.42%  '
$ws.Range("E45").Value = '  +3.51%  '
$ws.Range("D46").Value = '''40.66'
$ws.Range("E46").Value = '  +0.45%  '
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").Value = '''156.30'
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("E49").Value = '  +1.97%  '
$ws.Range("E50").Value = '  +0.75%  '
$ws.Range("D51").Value = '''21.92'
$ws.Range("E51").Value = '  -0.59%  '
